{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\n// ---------------------------------------------------------------------\n// 1) \"NOTE: ...\" paragraph: \"...downloaded from Kaggle(Link in README...\"\n//    becomes \"...downloaded from Kaggle (Link in README...\" (a single\n//    space is inserted between \"Kaggle\" and \"(Link\").\n// ---------------------------------------------------------------------\nconst body = context.document.body;\n\nconst kaggleHits = body.search(\"Kaggle(\", { matchCase: true });\nkaggleHits.load(\"items\");\nawait context.sync();\n\nif (kaggleHits.items.length === 0) {\n  throw new Error(\"Could not find 'Kaggle(' text to update.\");\n}\n\n// Replace \"Kaggle(\" with \"Kaggle (\" (adds the missing space).\nkaggleHits.items[0].insertText(\"Kaggle (\", \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Insert a brand-new paragraph (with the \"Update: - ...\" text) right\n//    after the blank paragraph that follows the \"NOTE:\" paragraph.\n// ---------------------------------------------------------------------\nconst noteHits = body.search(\"NOTE: - Once dataset is downloaded from\", { matchCase: true });\nnoteHits.load(\"items\");\nawait context.sync();\n\nif (noteHits.items.length === 0) {\n  throw new Error(\"Could not find the 'NOTE:' paragraph.\");\n}\n\nconst noteParagraph = noteHits.items[0].paragraphs.getFirst();\nconst blankParagraph = noteParagraph.getNext();\nblankParagraph.load(\"text\");\nawait context.sync();\n\nconst updateText =\n  \"Update: - Issue_and_varification folder is added with snapshot. Reference_image1 and \" +\n  \"Reference\" +\n  \"_image2 are for reference to make sure code is working fine. One issue while training was system detects an empty file while training so make sure to search for empty file in dataset folders and remove it forcefully.\";\n\nconst newParagraph = blankParagraph.insertParagraph(updateText, \"After\");\nnewParagraph.leftIndent = 18; // matches <w:ind w:left=\"360\"/> (360 twips = 18pt)\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) \"NOTE: ...\" paragraph: \"...downloaded from Kaggle(Link in README...\"\n#    becomes \"...downloaded from Kaggle (Link in README...\" (a single\n#    space is inserted between \"Kaggle\" and \"(Link\").\n# ---------------------------------------------------------------------\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Kaggle(\", $false, $false, $false, $false, $false, $true, 1, $false, \"Kaggle (\", 2)\nif (-not $found) {\n    throw \"Could not find 'Kaggle(' text to update.\"\n}\n\n# ---------------------------------------------------------------------\n# 2) Insert a brand-new paragraph (with the \"Update: - ...\" text) right\n#    after the blank paragraph that follows the \"NOTE:\" paragraph.\n# ---------------------------------------------------------------------\n$noteIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"NOTE: - Once dataset is downloaded from*\") {\n        $noteIndex = $i\n    }\n}\nif ($noteIndex -eq 0) {\n    throw \"Could not find the 'NOTE:' paragraph.\"\n}\n\n$blankPara = $d.Paragraphs.Item($noteIndex + 1)\n$blankPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($noteIndex + 2)\n$updateText = \"Update: - Issue_and_varification folder is added with snapshot. Reference_image1 and Reference_image2 are for reference to make sure code is working fine. One issue while training was system detects an empty file while training so make sure to search for empty file in dataset folders and remove it forcefully.\"\n$newPara.Range.Text = $updateText\n$newPara.LeftIndent = 18\n"}
